$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row: rename/retarget the columns and add the new kraken2 column.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "file"
$ws.Range("B1").Value = "reads aligned (minimap2)"
$ws.Range("C1").Value = "percent of total reads "
$ws.Range("D1").Value = "reads aligned (kraken2) "

# ---------------------------------------------------------------------------
# Species rows: the "k12" row is dropped and every row below it moves up
# one slot (rows 7-12 -> rows 6-11), while rows further down the sheet
# (the totals block starting at row 16) stay exactly where they were.
# Do this with plain cell writes rather than a full-row delete so the
# blank rows 12-15 aren't shifted.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "bsubtilis"
$ws.Range("B2").Value = 901704
$ws.Range("D2").Value = 636000

$ws.Range("A3").Value = "cneoformans"
$ws.Range("B3").Value = 116935
$ws.Range("D3").Value = 0

$ws.Range("A4").Value = "ecoli"
$ws.Range("B4").Value = 531245
$ws.Range("D4").Value = 187000

$ws.Range("A5").Value = "efaecalis"
$ws.Range("B5").Value = 649559
$ws.Range("D5").Value = 392000

$ws.Range("A6").Value = "lfermentum"
$ws.Range("B6").Value = 777255
$ws.Range("D6").Value = 514000

$ws.Range("A7").Value = "lmonocytogenes"
$ws.Range("B7").Value = 738425
$ws.Range("D7").Value = 471000

$ws.Range("A8").Value = "paeruginosa"
$ws.Range("B8").Value = 344018
$ws.Range("D8").Value = 184000

$ws.Range("A9").Value = "saureus"
$ws.Range("B9").Value = 702244
$ws.Range("D9").Value = 399000

$ws.Range("A10").Value = "scerevisiae"
$ws.Range("B10").Value = 107738
$ws.Range("D10").Value = 0

$ws.Range("A11").Value = "senterica"
$ws.Range("B11").Value = 532158
$ws.Range("D11").Value = 214000

# New column D values use the same "comma" number format as column B.
$ws.Range("D2:D11").NumberFormat = $ws.Range("B2").NumberFormat

# Row 12 used to hold "senterica_pb.ba" - it's no longer part of the table,
# so clear it out entirely (values + formula).
$ws.Range("A12:D12").Clear()

# C2:C11 percentage formulas shift reference-wise with the row, recompute
# them explicitly against the (unchanged) A19 total so they match the
# newly-relocated B column values.
$ws.Range("C2").Formula = "= 100 * (B2 / A19)"
$ws.Range("C3").Formula = "= 100 * (B3 / A19)"
$ws.Range("C4").Formula = "= 100 * (B4 / A19)"
$ws.Range("C5").Formula = "= 100 * (B5 / A19)"
$ws.Range("C6").Formula = "= 100 * (B6 / A19)"
$ws.Range("C7").Formula = "= 100 * (B7 / A19)"
$ws.Range("C8").Formula = "= 100 * (B8 / A19)"
$ws.Range("C9").Formula = "= 100 * (B9 / A19)"
$ws.Range("C10").Formula = "= 100 * (B10 / A19)"
$ws.Range("C11").Formula = "= 100 * (B11 / A19)"

# ---------------------------------------------------------------------------
# The grand-total SUM now only spans the 10 remaining species rows.
# ---------------------------------------------------------------------------
$ws.Range("A22").Formula = "=SUM(B2:B11)"

# ---------------------------------------------------------------------------
# Column widths - widen B and D for the new, longer headers. (The COM layer
# quantizes ColumnWidth to whole pixels before storing the OOXML "width" in
# characters, so the input here is chosen so it round-trips to the target
# stored widths of 24.5 / 21.5.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 23.67
$ws.Columns.Item(4).ColumnWidth = 20.67

# ---------------------------------------------------------------------------
# Selection moves to E1.
# ---------------------------------------------------------------------------
$ws.Range("E1").Select() | Out-Null
